$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format temporarily so numeric-looking values
# (e.g. "241.07", "1.000") are stored as literal text, matching the
# inlineStr cells produced by the source data feed.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.223.72"
$ws.Range("E2").Value = "  -2.32%  "

$ws.Range("D3").Value = "1.852.00"
$ws.Range("E3").Value = "  -1.25%  "

$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").Value = "0.6947"
$ws.Range("E5").Value = "  -5.97%  "

$ws.Range("D6").Value = "238.39"
$ws.Range("E6").Value = "  -1.76%  "

$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("D8").Value = "0.3065"
$ws.Range("E8").Value = "  -2.94%  "

$ws.Range("D9").Value = "0.07583"
$ws.Range("E9").Value = "  +5.31%  "

$ws.Range("D10").Value = "23.54"
$ws.Range("E10").Value = "  -4.30%  "

$ws.Range("D11").Value = "0.08087"
$ws.Range("E11").Value = "  -3.09%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.850.65"
$ws.Range("E12").Value = "  -2.92%  "

$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "0.7232"
$ws.Range("E13").Value = "  -3.51%  "

$ws.Range("D14").Value = "5.189"
$ws.Range("E14").Value = "  -3.72%  "

$ws.Range("D15").Value = "89.03"
$ws.Range("E15").Value = "  -3.55%  "

$ws.Range("D16").Value = "29.211.00"
$ws.Range("E16").Value = "  -2.39%  "

$ws.Range("D17").Value = "5.772"
$ws.Range("E17").Value = "  -5.24%  "

$ws.Range("D18").Value = "241.07"

$ws.Range("D19").Value = "0.000007723"
$ws.Range("E19").Value = "  -1.36%  "

$ws.Range("D20").Value = "13.08"
$ws.Range("E20").Value = "  -3.45%  "

$ws.Range("D21").Value = "0.9997"
$ws.Range("E21").Value = "  -0.11%  "

$ws.Range("D22").Value = "2.093.86"
$ws.Range("E22").Value = "  -1.98%  "

$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  -0.10%  "

$ws.Range("D24").Value = "7.606"
$ws.Range("E24").Value = "  -5.14%  "

$ws.Range("D25").Value = "9.020"
$ws.Range("E25").Value = "  -2.72%  "

$ws.Range("D26").Value = "161.48"

$ws.Range("D27").Value = "0.1454"
$ws.Range("E27").Value = "  -6.33%  "

$ws.Range("D28").Value = "18.04"
$ws.Range("E28").Value = "  -3.30%  "

$ws.Range("D29").Value = "1.933"
$ws.Range("E29").Value = "  -4.71%  "

$ws.Range("D30").Value = "1.396"
$ws.Range("E30").Value = "  -7.46%  "

$ws.Range("D31").Value = "1.502"
$ws.Range("E31").Value = "  -2.16%  "

$ws.Range("D32").Value = "4.429"
$ws.Range("E32").Value = "  -3.60%  "

$ws.Range("D33").Value = "4.047"
$ws.Range("E33").Value = "  -5.19%  "

$ws.Range("D34").Value = "0.05227"
$ws.Range("E34").Value = "  -1.66%  "

$ws.Range("D35").Value = "1.192"
$ws.Range("E35").Value = "  -3.38%  "

$ws.Range("D36").Value = "0.7073"
$ws.Range("E36").Value = "  -5.49%  "

$ws.Range("D37").Value = "1.001"
$ws.Range("E37").Value = "  +0.16%  "

$ws.Range("D38").Value = "2.666"
$ws.Range("E38").Value = "  -1.17%  "

$ws.Range("D39").Value = "0.01861"
$ws.Range("E39").Value = "  -5.27%  "

$ws.Range("D40").Value = "2.695"
$ws.Range("E40").Value = "  -2.24%  "

$ws.Range("D41").Value = "0.9123"
$ws.Range("E41").Value = "  +5.91%  "

$ws.Range("D42").Value = "5.960"
$ws.Range("E42").Value = "  -2.63%  "

$ws.Range("D43").Value = "0.4292"
$ws.Range("E43").Value = "  -5.52%  "

$ws.Range("D44").Value = "1.044.90"
$ws.Range("E44").Value = "  -5.40%  "

$ws.Range("D45").Value = "69.49"
$ws.Range("E45").Value = "  -3.95%  "

$ws.Range("D46").Value = "1.000"
$ws.Range("E46").Value = "  -0.20%  "

$ws.Range("D47").Value = "102.34"
$ws.Range("E47").Value = "  -1.85%  "

$ws.Range("D48").Value = "7.236"
$ws.Range("E48").Value = "  -4.85%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "9.282"
$ws.Range("E49").Value = "  -2.45%  "

$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "1.738"
$ws.Range("E50").Value = "  -6.32%  "

$ws.Range("D51").Value = "1.993.18"
$ws.Range("E51").Value = "  -2.20%  "

# Restore the default cell style now that the text values are set,
# so the cells do not retain an explicit text number format.
$ws.Range("D2:D51").Style = "Normal"
